$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 40, shifting existing rows 40:154 down to 41:154
$ws.Rows.Item(40).Insert()

# Populate the new row 40 with the new data entry
$ws.Range("A40").Value = 5
$ws.Range("B40").Value = "Macroferia Regional de Talca"
$ws.Range("C40").Value = "Maule"
$ws.Range("D40").Value = 45238
$ws.Range("D40").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E40").Value = 7
$ws.Range("F40").Value = 100112022
$ws.Range("G40").Value = "Arveja Verde"
$ws.Range("H40").Value = "Sin especificar"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 500
$ws.Range("K40").Value = 20000
$ws.Range("L40").Value = 20000
$ws.Range("M40").Value = 20000
$ws.Range("N40").Value = "`$/saco 25 kilos"
$ws.Range("O40").Value = "Región del Maule"
$ws.Range("P40").Value = 800
$ws.Range("Q40").Value = 25
$ws.Range("R40").Value = "Hortaliza"
